$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 137410
$ws.Range("J63").Value = 137410
$ws.Range("L63").Value = 137410
$ws.Range("N63").Value = -138658

$ws.Range("H66").Value = 137410
$ws.Range("J66").Value = 137410
$ws.Range("L66").Value = 412230
$ws.Range("N66").Value = -418470

$ws.Range("H75").Value = 42000
$ws.Range("J75").Value = 42000
$ws.Range("L75").Value = 42000
$ws.Range("N75").Value = -43872

$ws.Range("H78").Value = 42000
$ws.Range("J78").Value = 42000
$ws.Range("L78").Value = 126000
$ws.Range("N78").Value = -135360

$ws.Range("H93").Value = 42000
$ws.Range("J93").Value = 42000
$ws.Range("L93").Value = 42000
$ws.Range("N93").Value = -46992

$ws.Range("H129").Value = 1388.0278
$ws.Range("I129").Value = 632.7143
$ws.Range("J129").Value = 1570.3448
$ws.Range("K129").Value = 1898.1429
$ws.Range("L129").Value = 4711.0344
$ws.Range("M129").Value = 3101.8571
$ws.Range("N129").Value = -14711.0344

$ws.Range("H130").Value = 51802.5
$ws.Range("J130").Value = 51802.5
$ws.Range("L130").Value = 51802.5
$ws.Range("N130").Value = -61842.5

$ws.Range("H137").Value = 1936.4231
$ws.Range("I137").Value = 1472.35
$ws.Range("J137").Value = 3483.3333
$ws.Range("K137").Value = 4417.049999999999
$ws.Range("L137").Value = 10449.9999
$ws.Range("M137").Value = -1867.049999999999
$ws.Range("N137").Value = -15549.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32804.15
$ws.Range("I32").Value = 6250.7144
$ws.Range("K32").Value = 6250.7144
$ws.Range("M32").Value = -5963.7144

$ws.Range("H55").Value = 20740.75
$ws.Range("J55").Value = 23696.857
$ws.Range("L55").Value = 23696.857
$ws.Range("N55").Value = -24326.857

$ws.Range("H80").Value = 26677.111
$ws.Range("J80").Value = 26677.111
$ws.Range("L80").Value = 26677.111
$ws.Range("N80").Value = -28673.111

$ws.Range("H83").Value = 26677.111
$ws.Range("J83").Value = 26677.111
$ws.Range("L83").Value = 80031.333
$ws.Range("N83").Value = -90015.333

$ws.Range("H103").Value = 35951.25
$ws.Range("J103").Value = 35951.25
$ws.Range("L103").Value = 35951.25
$ws.Range("N103").Value = -38295.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18898.166
$ws.Range("I82").Value = 10257
$ws.Range("J82").Value = 20626.4
$ws.Range("K82").Value = 10257
$ws.Range("L82").Value = 20626.4
$ws.Range("M82").Value = -9874
$ws.Range("N82").Value = -21392.4

$ws.Range("H85").Value = 18898.166
$ws.Range("I85").Value = 10257
$ws.Range("J85").Value = 20626.4
$ws.Range("K85").Value = 10257
$ws.Range("L85").Value = 20626.4
$ws.Range("M85").Value = -8931
$ws.Range("N85").Value = -23278.4

$ws.Range("H99").Value = 3412.5
$ws.Range("I99").Value = 1960
$ws.Range("K99").Value = 1960
$ws.Range("M99").Value = -462

$ws.Range("H122").Value = 50914.285
$ws.Range("J122").Value = 50914.285
$ws.Range("L122").Value = 50914.285
$ws.Range("N122").Value = -60714.285

$ws.Range("H135").Value = 53780
$ws.Range("J135").Value = 53780
$ws.Range("L135").Value = 53780
$ws.Range("N135").Value = -63920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1882.3489
$ws.Range("I31").Value = 1346.931
$ws.Range("J31").Value = 2991.4285
$ws.Range("K31").Value = 1346.931
$ws.Range("L31").Value = 2991.4285
$ws.Range("M31").Value = -1051.931
$ws.Range("N31").Value = -3581.4285

$ws.Range("H34").Value = 1882.3489
$ws.Range("I34").Value = 1346.931
$ws.Range("J34").Value = 2991.4285
$ws.Range("K34").Value = 1346.931
$ws.Range("L34").Value = 2991.4285
$ws.Range("M34").Value = -1144.931
$ws.Range("N34").Value = -3395.4285

$ws.Range("H41").Value = 14598.429
$ws.Range("J41").Value = 18246
$ws.Range("L41").Value = 18246
$ws.Range("N41").Value = -19102

$ws.Range("H50").Value = 8592.333
$ws.Range("J50").Value = 9041.375
$ws.Range("L50").Value = 9041.375
$ws.Range("N50").Value = -10291.375

$ws.Range("H51").Value = 7545
$ws.Range("J51").Value = 10590
$ws.Range("L51").Value = 10590
$ws.Range("N51").Value = -12062

$ws.Range("H60").Value = 14026.875
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 14026.875
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 14026.875
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -15048.875

$ws.Range("H61").Value = 7545
$ws.Range("J61").Value = 10590
$ws.Range("L61").Value = 10590
$ws.Range("N61").Value = -11286

$ws.Range("H68").Value = 34014.5
$ws.Range("J68").Value = 35263.332
$ws.Range("L68").Value = 35263.332
$ws.Range("N68").Value = -36761.332

$ws.Range("H71").Value = 34014.5
$ws.Range("J71").Value = 35263.332
$ws.Range("L71").Value = 105789.996
$ws.Range("N71").Value = -113277.996

$ws.Range("H109").Value = 20694
$ws.Range("J109").Value = 20694
$ws.Range("L109").Value = 20694
$ws.Range("N109").Value = -22774

$ws.Range("H127").Value = 50655
$ws.Range("J127").Value = 50655
$ws.Range("L127").Value = 50655
$ws.Range("N127").Value = -60575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 787.3333
$ws.Range("J97").Value = 787.3333
$ws.Range("L97").Value = 2361.9999
$ws.Range("N97").Value = -3353.9999

$ws.Range("H107").Value = 906.4091
$ws.Range("I107").Value = 629.0769
$ws.Range("J107").Value = 1307
$ws.Range("K107").Value = 1887.2307
$ws.Range("L107").Value = 3921
$ws.Range("M107").Value = 32.76929999999993
$ws.Range("N107").Value = -7761

$ws.Range("H137").Value = 1643.56
$ws.Range("I137").Value = 981.7059
$ws.Range("J137").Value = 3050
$ws.Range("K137").Value = 2945.1177
$ws.Range("L137").Value = 9150
$ws.Range("M137").Value = 2154.8823
$ws.Range("N137").Value = -19350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 17822.7
$ws.Range("J57").Value = 17822.7
$ws.Range("L57").Value = 17822.7
$ws.Range("N57").Value = -19462.7

$ws.Range("H93").Value = 9817.1
$ws.Range("J93").Value = 9817.1
$ws.Range("L93").Value = 9817.1
$ws.Range("N93").Value = -13561.1

$ws.Range("H102").Value = 2059.9
$ws.Range("I102").Value = 1870.4706
$ws.Range("J102").Value = 3133.3333
$ws.Range("K102").Value = 1870.4706
$ws.Range("L102").Value = 3133.3333
$ws.Range("M102").Value = -248.4706000000001
$ws.Range("N102").Value = -6377.3333

$ws.Range("H122").Value = 2254.3157
$ws.Range("I122").Value = 2254.9333
$ws.Range("J122").Value = 2252
$ws.Range("K122").Value = 6764.7999
$ws.Range("L122").Value = 6756
$ws.Range("M122").Value = -4314.7999
$ws.Range("N122").Value = -11656

$ws.Range("H123").Value = 21730.4
$ws.Range("J123").Value = 21730.4
$ws.Range("L123").Value = 21730.4
$ws.Range("N123").Value = -26630.4

$ws.Range("H124").Value = 48692
$ws.Range("J124").Value = 48692
$ws.Range("L124").Value = 48692
$ws.Range("N124").Value = -58512

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1859.8
$ws.Range("I68").Value = 1220.4
$ws.Range("J68").Value = 2499.2
$ws.Range("K68").Value = 1220.4
$ws.Range("L68").Value = 2499.2
$ws.Range("M68").Value = -471.4000000000001
$ws.Range("N68").Value = -3997.2

$ws.Range("H71").Value = 1859.8
$ws.Range("I71").Value = 1220.4
$ws.Range("J71").Value = 2499.2
$ws.Range("K71").Value = 6102
$ws.Range("L71").Value = 12496
$ws.Range("M71").Value = -2358
$ws.Range("N71").Value = -19984

$ws.Range("H125").Value = 43683.332
$ws.Range("J125").Value = 43683.332
$ws.Range("L125").Value = 43683.332
$ws.Range("N125").Value = -53523.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws.Range("H122").Value = 2563.3513
$ws.Range("I122").Value = 2053.2593
$ws.Range("J122").Value = 3940.6
$ws.Range("K122").Value = 6159.777900000001
$ws.Range("L122").Value = 11821.8
$ws.Range("M122").Value = -3709.777900000001
$ws.Range("N122").Value = -16721.8
